$wb = $excel.ActiveWorkbook

# Swap the values in row 2 of KonnectionsUsers: (Vijay,Sunil,Vijay,Sunil) -> (Sunil,Vijay,Sunil,Vijay)
$wsUsers = $wb.Worksheets.Item("KonnectionsUsers")
$wsUsers.Range("A2").Value = "Sunil"
$wsUsers.Range("B2").Value = "Vijay"
$wsUsers.Range("C2").Value = "Sunil"
$wsUsers.Range("D2").Value = "Vijay"

# Update selections on both affected sheets before changing the active tab.
$wsLogins = $wb.Worksheets.Item("Logins")
$wsLogins.Range("B5").Select()

$wsUsers.Range("C8").Select()

# Make KonnectionsUsers the active sheet (was Logins).
$wsUsers.Activate()
